$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.316.38"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.039.38"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0799"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").Value = "2.345.21"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.748"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "2.046.36"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "37.199.84"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.126"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.117"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0219"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.81%  "
$ws.Range("D41").Value = "1.492.33"
$ws.Range("E41").Value = "  +2.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0940"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("E46").Value = "  -5.87%  "
$ws.Range("E47").Value = "  -4.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.20%  "
$ws.Range("D51").Value = "2.231.79"
$ws.Range("E51").Value = "  -2.16%  "
